$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update the two background hitcount values from 5-day to 1-day counts
$ws.Range("F5").Value = 56
$ws.Range("F8").Value = 55

# Reflect the active selection recorded in the saved file
$ws.Range("G13").Select()
